$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = "58233443"
$ws.Range("C2").Value = "FCT910932700171337728"
$ws.Range("F2").Value = "171.90"
$ws.Range("I2").Value = "58233443+1"
